$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.07636889456857432
$ws.Range("D2").Value = 0.1364922641372459
$ws.Range("E2").Value = 0.1522011239481671
$ws.Range("F2").Value = 1.923069012955807
$ws.Range("G2").Value = 1.265832160430094
$ws.Range("H2").Value = 1.190944251512263
$ws.Range("I2").Value = 0.5818441923732784
$ws.Range("J2").Value = 0.1994494840745205
$ws.Range("K2").Value = 1.17305898848295

$ws.Range("B3").Value = 0.06694771805051403
$ws.Range("D3").Value = 0.1328050409581749
$ws.Range("E3").Value = 0.1486666221905821
$ws.Range("F3").Value = 1.917959767128522
$ws.Range("G3").Value = 1.261498715295843
$ws.Range("H3").Value = 1.195280968651019
$ws.Range("I3").Value = 0.5908596728898825
$ws.Range("J3").Value = 0.195242944255547
$ws.Range("K3").Value = 1.062811734370115

$ws.Range("B4").Value = 0.0611510855129751
$ws.Range("D4").Value = 0.1305856332785709
$ws.Range("E4").Value = 0.1465722913638423
$ws.Range("F4").Value = 1.916055699280406
$ws.Range("G4").Value = 1.259788548364327
$ws.Range("H4").Value = 1.198610323401184
$ws.Range("I4").Value = 0.5967222617122561
$ws.Range("J4").Value = 0.1927804758424188
$ws.Range("K4").Value = 0.9953671376197235

$ws.Range("B5").Value = 0.05878608492604087
$ws.Range("D5").Value = 0.129692472417986
$ws.Range("E5").Value = 0.1457379132161165
$ws.Range("F5").Value = 1.915589205066524
$ws.Range("G5").Value = 1.259329902759191
$ws.Range("H5").Value = 1.200134463004389
$ws.Range("I5").Value = 0.5991934815124802
$ws.Range("J5").Value = 0.1918072019902795
$ws.Range("K5").Value = 0.9679456649982683

$ws.Range("B6").Value = 0.05839321332830139
$ws.Range("D6").Value = 0.1295448459535322
$ws.Range("E6").Value = 0.1456005177214479
$ws.Range("F6").Value = 1.915530414572189
$ws.Range("G6").Value = 1.259268112893679
$ws.Range("H6").Value = 1.200397649421461
$ws.Range("I6").Value = 0.599608786287801
$ws.Range("J6").Value = 0.1916474128741612
$ws.Range("K6").Value = 0.9633961545734167

$ws.Range("B7").Value = 0.06111920147648675
$ws.Range("D7").Value = 0.1305735421032352
$ws.Range("E7").Value = 0.1465609613968049
$ws.Range("F7").Value = 1.916048155948502
$ws.Range("G7").Value = 1.259781399225716
$ws.Range("H7").Value = 1.198630200981327
$ws.Range("I7").Value = 0.5967552568512957
$ws.Range("J7").Value = 0.1927672277275363
$ws.Range("K7").Value = 0.9949970673461621

$ws.Range("B8").Value = 0.07312307225809889
$ws.Range("D8").Value = 0.1352116952134566
$ws.Range("E8").Value = 0.1509666804501713
$ws.Range("F8").Value = 1.921051044444752
$ws.Range("G8").Value = 1.264140224478638
$ws.Range("H8").Value = 1.192301065510719
$ws.Range("I8").Value = 0.5848848093858017
$ws.Range("J8").Value = 0.1979740467339823
$ws.Range("K8").Value = 1.134994574914515

$ws.Range("B9").Value = 0.09656033088629101
$ws.Range("D9").Value = 0.1446586314353908
$ws.Range("E9").Value = 0.1602088738963729
$ws.Range("F9").Value = 1.940676929329499
$ws.Range("G9").Value = 1.280267582551019
$ws.Range("H9").Value = 1.18518973895803
$ws.Range("I9").Value = 0.5642047440643019
$ws.Range("J9").Value = 0.2091433659286679
$ws.Range("K9").Value = 1.411488946557938

$ws.Range("B10").Value = 0.1137094744520084
$ws.Range("D10").Value = 0.1518117268596342
$ws.Range("E10").Value = 0.1673682902234361
$ws.Range("F10").Value = 1.961128556565043
$ws.Range("G10").Value = 1.296792252341874
$ws.Range("H10").Value = 1.183212859717884
$ws.Range("I10").Value = 0.5505982902781046
$ws.Range("J10").Value = 0.2179402966738877
$ws.Range("K10").Value = 1.615839414146706

$ws.Range("B11").Value = 0.1214942882633636
$ws.Range("D11").Value = 0.1551116547513089
$ws.Range("E11").Value = 0.1707059209507449
$ws.Range("F11").Value = 1.971753351131113
$ws.Range("G11").Value = 1.305337222234215
$ws.Range("H11").Value = 1.183022515343879
$ws.Range("I11").Value = 0.5447537428457956
$ws.Range("J11").Value = 0.2220719945827483
$ws.Range("K11").Value = 1.709071488404732

$ws.Range("B12").Value = 0.1244396708914763
$ws.Range("D12").Value = 0.1563678143911034
$ws.Range("E12").Value = 0.1719814294398887
$ws.Range("F12").Value = 1.97596748181013
$ws.Range("G12").Value = 1.308721703115538
$ws.Range("H12").Value = 1.183052670532675
$ws.Range("I12").Value = 0.5425902675767151
$ws.Range("J12").Value = 0.2236553398935541
$ws.Range("K12").Value = 1.744415116689652

$ws.Range("B13").Value = 0.1238054468942948
$ws.Range("D13").Value = 0.1560969876153422
$ws.Range("E13").Value = 0.171706208973049
$ws.Range("F13").Value = 1.97505139750595
$ws.Range("G13").Value = 1.307986167555782
$ws.Range("H13").Value = 1.183041624536571
$ws.Range("I13").Value = 0.5430539973644368
$ws.Range("J13").Value = 0.2233135025111608
$ws.Range("K13").Value = 1.73680152128054

$ws.Range("B14").Value = 0.1217366590385467
$ws.Range("D14").Value = 0.1552148688522124
$ws.Range("E14").Value = 0.1708106248868475
$ws.Range("F14").Value = 1.972096222890855
$ws.Range("G14").Value = 1.305612680518976
$ws.Range("H14").Value = 1.18302294563486
$ws.Range("I14").Value = 0.5445747552153577
$ws.Range("J14").Value = 0.2222018809101769
$ws.Range("K14").Value = 1.71197845757905

$ws.Range("B15").Value = 0.1204691278912691
$ws.Range("D15").Value = 0.1546753967238033
$ws.Range("E15").Value = 0.1702635672147039
$ws.Range("F15").Value = 1.970310957188559
$ws.Range("G15").Value = 1.304178240935187
$ws.Range("H15").Value = 1.183024826831371
$ws.Range("I15").Value = 0.5455127424451192
$ws.Range("J15").Value = 0.2215234261256569
$ws.Range("K15").Value = 1.696778636615591

$ws.Range("B16").Value = 0.113200371348114
$ws.Range("D16").Value = 0.1515969883706987
$ws.Range("E16").Value = 0.1671517940103371
$ws.Range("F16").Value = 1.960460852766644
$ws.Range("G16").Value = 1.296254576507096
$ws.Range("H16").Value = 1.183239589617074
$ws.Range("I16").Value = 0.5509871982612751
$ws.Range("J16").Value = 0.2176729015987888
$ws.Range("K16").Value = 1.609751913563287

$ws.Range("B17").Value = 0.1087368792024961
$ws.Range("D17").Value = 0.1497202110268887
$ws.Range("E17").Value = 0.1652635137437457
$ws.Range("F17").Value = 1.954757073103693
$ws.Range("G17").Value = 1.291657543254701
$ws.Range("H17").Value = 1.18355312257583
$ws.Range("I17").Value = 0.5544340605913156
$ws.Range("J17").Value = 0.2153440623163476
$ws.Range("K17").Value = 1.556433181494185

$ws.Range("B18").Value = 0.1061680628951933
$ws.Range("D18").Value = 0.1486450667265018
$ws.Range("E18").Value = 0.1641850276040131
$ws.Range("F18").Value = 1.951600726265013
$ws.Range("G18").Value = 1.289110140789603
$ws.Range("H18").Value = 1.183800170551706
$ws.Range("I18").Value = 0.5564490805536515
$ws.Range("J18").Value = 0.2140168005054761
$ws.Range("K18").Value = 1.525791286673098

$ws.Range("B19").Value = 0.1052980484439416
$ws.Range("D19").Value = 0.1482817865521895
$ws.Range("E19").Value = 0.1638211763752224
$ws.Range("F19").Value = 1.950553369903716
$ws.Range("G19").Value = 1.288264213608358
$ws.Range("H19").Value = 1.183895265866539
$ws.Range("I19").Value = 0.5571369063106815
$ws.Range("J19").Value = 0.213569509772924
$ws.Range("K19").Value = 1.515420883588831

$ws.Range("B20").Value = 0.1092121857205797
$ws.Range("D20").Value = 0.1499195496356691
$ws.Range("E20").Value = 0.1654637375166814
$ws.Range("F20").Value = 1.955351378474887
$ws.Range("G20").Value = 1.292136891905415
$ws.Range("H20").Value = 1.183512839948662
$ws.Range("I20").Value = 0.5540637743382817
$ws.Range("J20").Value = 0.2155907054413149
$ws.Range("K20").Value = 1.562106402302277

$ws.Range("B21").Value = 0.1223443831441244
$ws.Range("D21").Value = 0.1554737912535558
$ws.Range("E21").Value = 0.1710733642335001
$ws.Range("F21").Value = 1.97295904629685
$ws.Range("G21").Value = 1.306305788796976
$ws.Range("H21").Value = 1.183025655150942
$ws.Range("I21").Value = 0.5441267213392607
$ws.Range("J21").Value = 0.2225278813324394
$ws.Range("K21").Value = 1.719268547240631

$ws.Range("B22").Value = 0.1309120369352854
$ws.Range("D22").Value = 0.1591419526185689
$ws.Range("E22").Value = 0.1748073135238144
$ws.Range("F22").Value = 1.985578877215701
$ws.Range("G22").Value = 1.316433054372652
$ws.Range("H22").Value = 1.183303288629617
$ws.Range("I22").Value = 0.5379221840467405
$ws.Range("J22").Value = 0.2271711167427668
$ws.Range("K22").Value = 1.822208269578425

$ws.Range("B23").Value = 0.1263407554727678
$ws.Range("D23").Value = 0.1571807149426832
$ws.Range("E23").Value = 0.1728082351168538
$ws.Range("F23").Value = 1.978741417414085
$ws.Range("G23").Value = 1.310948310473464
$ws.Range("H23").Value = 1.18310047684173
$ws.Range("I23").Value = 0.541207102035294
$ws.Range("J23").Value = 0.2246828996468366
$ws.Range("K23").Value = 1.767246942621341

$ws.Range("B24").Value = 0.1089973080242714
$ws.Range("D24").Value = 0.149829416685975
$ws.Range("E24").Value = 0.1653731941994252
$ws.Range("F24").Value = 1.955082310463894
$ws.Range("G24").Value = 1.291919880996261
$ws.Range("H24").Value = 1.18353084368735
$ws.Range("I24").Value = 0.5542310768079748
$ws.Range("J24").Value = 0.2154791619002765
$ws.Range("K24").Value = 1.559541502530237

$ws.Range("B25").Value = 0.09023170191592556
$ws.Range("D25").Value = 0.142065547870871
$ws.Range("E25").Value = 0.1576439187250998
$ws.Range("F25").Value = 1.934311789187461
$ws.Range("G25").Value = 1.275087563672628
$ws.Range("H25").Value = 1.186544411090253
$ws.Range("I25").Value = 0.5695207754950804
$ws.Range("J25").Value = 0.206018474801084
$ws.Range("K25").Value = 1.336478582364293
